# "Add files via upload" - refreshed Gpmp coefficient row (row 23) with
# newly recomputed (much smaller magnitude) regression coefficients, and
# left the selection on the cell the author was last looking at (J24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 ("Gpmp") B:J - replace the old coefficients with the refreshed values.
$ws.Range("B23").Value = -0.0033281476303707201
$ws.Range("C23").Value = -0.0025440984127414398
$ws.Range("D23").Value = -0.0039656922651073404
$ws.Range("E23").Value = -0.0040022856049998197
$ws.Range("F23").Value = -0.0035005807957033202
$ws.Range("G23").Value = -0.0043192833730332303
$ws.Range("H23").Value = -0.0035229329807558202
$ws.Range("I23").Value = -0.0034524162627164799
$ws.Range("J23").Value = -0.0034388956141986199

# Leave the cursor on J24, matching the author's saved selection.
$ws.Range("J24").Select()
